# Leave card update: add 2023 year block (rows 74-77 leave entries) and
# extend the monthly PERIOD rows through row 132 (table grows from
# A8:K130 to A8:K132). Mirrors commit "Leave 3/9/2023 12:08 AM".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# 1) Grow Table1 by two rows (130 -> 132) using the table's own "add
#    row" operation so the table definition (ref, calculated columns)
#    stays consistent.
# ---------------------------------------------------------------------
$lo = $ws.ListObjects.Item(1)
$lo.ListRows.Add() | Out-Null
$lo.ListRows.Add() | Out-Null

# ---------------------------------------------------------------------
# 2) Re-create the row formatting around the new table rows.
#    Before the insert, row 130 held the special "last row" formatting
#    (thicker bottom border etc.). That row keeps its place; the two
#    brand new rows (131/132) come in unformatted, and the old
#    "last row" look needs to move onto the new last row (132).
# ---------------------------------------------------------------------

# 2a. Push the previous "last row" look from row 130 onto the new last
#     row, 132.
$ws.Range("A130:K130").Copy() | Out-Null
$ws.Range("A132:K132").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

# 2b. Give rows 130-131 the regular interior-row look (copied from row
#     129, which already had it).
$ws.Range("A129:K129").Copy() | Out-Null
$ws.Range("A130:K131").PasteSpecial(-4122) | Out-Null    # xlPasteFormats
$excel.CutCopyMode = 0

# 2c. Column A keeps the plain date style on every row, including the
#     new final row 132 (only columns B:K use the heavier "last row"
#     border there).
$ws.Range("A129").Copy() | Out-Null
$ws.Range("A132").PasteSpecial(-4122) | Out-Null         # xlPasteFormats
$excel.CutCopyMode = 0

# 2d. Re-apply the calculated-column formula text (PasteSpecial above
#     only moved formatting) for the two brand-new rows.
$balFormula = '=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"",Table1[[#This Row],[EARNED]])'
$ws.Range("G131").Formula = $balFormula
$ws.Range("G132").Formula = $balFormula

# ---------------------------------------------------------------------
# 3) New "2023" year-divider row (row 74): column A becomes a bold,
#    centered text label "2023" (matching the existing 2018/2019/2020/
#    2022 dividers already in the sheet), everything else on the row
#    is untouched.
# ---------------------------------------------------------------------
$ws.Range("A61").Copy() | Out-Null                        # the "2022" divider cell
$ws.Range("A74").PasteSpecial(-4122) | Out-Null           # xlPasteFormats
$excel.CutCopyMode = 0
$ws.Range("A74").Value2 = "'2023"                          # force text, not a number

# ---------------------------------------------------------------------
# 4) Leave entries for Jan 2023 (rows 75-76) and a new SPL row (77).
# ---------------------------------------------------------------------
$ws.Range("A75").Value2 = 44927
$ws.Range("B75").Value2 = "SL(1-0-0)"
$ws.Range("C75").Value2 = 1.25
$ws.Range("H75").Value2 = 1
$ws.Range("K75").Value2 = 44951

$ws.Range("A76").Value2 = 44958
$ws.Range("B76").Value2 = "SL(1-0-0)"
$ws.Range("C76").Value2 = 1.25
$ws.Range("H76").Value2 = 1
$ws.Range("K76").Value2 = 44980

$ws.Range("B77").Value2 = "SP(1-0-0)"
$ws.Range("K77").Value2 = 44979

# K column (dates) on rows 75-77 switch from the plain style to the
# short-date style already used on K71/K72.
$ws.Range("K71").Copy() | Out-Null
$ws.Range("K75:K77").PasteSpecial(-4122) | Out-Null        # xlPasteFormats
$excel.CutCopyMode = 0
# re-set the values (format-only paste does not touch them, but this
# keeps the script robust to engine paste-format quirks)
$ws.Range("K75").Value2 = 44951
$ws.Range("K76").Value2 = 44980
$ws.Range("K77").Value2 = 44979

# ---------------------------------------------------------------------
# 5) Column A monthly PERIOD dates for rows 78-132 (first of each
#    month, continuing the existing sequence).
# ---------------------------------------------------------------------
$periodDates = @{
  78 = 44986;  79 = 45017;  80 = 45047;  81 = 45078;  82 = 45108
  83 = 45139;  84 = 45170;  85 = 45200;  86 = 45231;  87 = 45261
  88 = 45292;  89 = 45323;  90 = 45352;  91 = 45383;  92 = 45413
  93 = 45444;  94 = 45474;  95 = 45505;  96 = 45536;  97 = 45566
  98 = 45597;  99 = 45627; 100 = 45658; 101 = 45689; 102 = 45717
 103 = 45748; 104 = 45778; 105 = 45809; 106 = 45839; 107 = 45870
 108 = 45901; 109 = 45931; 110 = 45962; 111 = 45992; 112 = 46023
 113 = 46054; 114 = 46082; 115 = 46113; 116 = 46143; 117 = 46174
 118 = 46204; 119 = 46235; 120 = 46266; 121 = 46296; 122 = 46327
 123 = 46357; 124 = 46388; 125 = 46419; 126 = 46447; 127 = 46478
 128 = 46508; 129 = 46539; 130 = 46569; 131 = 46600; 132 = 46631
}

foreach ($r in $periodDates.Keys) {
  $ws.Cells.Item($r, 1).Value2 = $periodDates[$r]
}

# ---------------------------------------------------------------------
# 6) Footer: certifying officer changed.
# ---------------------------------------------------------------------
$ps = $ws.PageSetup
$ps.OddFooter = "&L" + "`n" + "PREPARED BY: ___________________" + "`n" + "DATE: &D, &T&C" + "`n" + "CERTIFIED CORRECT BY: NANETTE B. SUSA" + "`n" + "                                              OIC - HRMO&RPage &P of &N"

# ---------------------------------------------------------------------
# 7) Leave the final selection close to the sheet's new scroll target.
# ---------------------------------------------------------------------
$ws.Range("E5").Select() | Out-Null
$ws.Range("B81").Activate() | Out-Null

Write-Host "edit complete"
